$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: Apply header style to newly-needed header cells (D1:O1) ----
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1:O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---- Step 2: Rewrite header row (D1:O1) with the new column names ----
$ws.Range("D1").Value = "CV Train F1"
$ws.Range("E1").Value = "CV Test F1"
$ws.Range("F1").Value = "Validation F1"
$ws.Range("G1").Value = "CV Train Precision"
$ws.Range("H1").Value = "CV Test Precision"
$ws.Range("I1").Value = "Validation Precision"
$ws.Range("J1").Value = "CV Train Recall"
$ws.Range("K1").Value = "CV Test Recall"
$ws.Range("L1").Value = "Validation Recall"
$ws.Range("M1").Value = "Y Val (Validation)"
$ws.Range("N1").Value = "Y Pred (Validation)"
$ws.Range("O1").Value = "Seed"

# ---- Step 3: Rewrite data rows 2-6 with the re-computed 70/30 split metrics ----
# Row 2
$ws.Range("A2").Value = @'
Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector', RandomUnderSampler(random_state=42)),
                ('model',
                 BaggingClassifier(estimator=DecisionTreeClassifier(criterion='entropy',
                                                                    max_depth=5,
                                                                    max_features='log2',
                                                                    min_samples_leaf=6,
                                                                    min_samples_split=5,
                                                                    random_state=42),
                                   n_estimators=5, random_state=42))])
'@
$ws.Range("B2").Value = 0.7366666666666666
$ws.Range("C2").Value = @'
{'selector': RandomUnderSampler(random_state=42), 'scaler': MinMaxScaler(), 'model__n_estimators': 5, 'model__estimator__min_samples_split': 5, 'model__estimator__min_samples_leaf': 6, 'model__estimator__max_features': 'log2', 'model__estimator__max_depth': 5, 'model__estimator__criterion': 'entropy', 'model__estimator__class_weight': None}
'@
$ws.Range("D2").Value = 0.8042932842433356
$ws.Range("E2").Value = 0.4381470307470308
$ws.Range("F2").Value = 0.7307692307692307
$ws.Range("G2").Value = 0.8870676688220692
$ws.Range("H2").Value = 0.5384761904761906
$ws.Range("I2").Value = 0.6333333333333333
$ws.Range("J2").Value = 0.7392000000000001
$ws.Range("K2").Value = 0.39232
$ws.Range("L2").Value = 0.8636363636363636
$ws.Range("M2").Value = @'
[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1 0 1 1 1 1 0 0 0 0 1 0 1]
'@
$ws.Range("N2").Value = @'
[1 1 1 1 1 1 1 0 1 1 0 1 1 1 0 1 0 1 1 1 1 1 1 1 1 0 1 1 1 1 1 1 0 1 1 1]
'@
$ws.Range("O2").Value = 42
$ws.Rows.Item(2).AutoFit()

# Row 3
$ws.Range("A3").Value = @'
Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector', RandomUnderSampler(random_state=42)),
                ('model',
                 BaggingClassifier(estimator=DecisionTreeClassifier(max_depth=4,
                                                                    max_features='log2',
                                                                    min_samples_leaf=5,
                                                                    min_samples_split=5,
                                                                    random_state=42),
                                   random_state=42))])
'@
$ws.Range("B3").Value = 0.7238095238095237
$ws.Range("C3").Value = @'
{'selector': RandomUnderSampler(random_state=42), 'scaler': MinMaxScaler(), 'model__n_estimators': 10, 'model__estimator__min_samples_split': 5, 'model__estimator__min_samples_leaf': 5, 'model__estimator__max_features': 'log2', 'model__estimator__max_depth': 4, 'model__estimator__criterion': 'gini', 'model__estimator__class_weight': None}
'@
$ws.Range("D3").Value = 0.8281337896127793
$ws.Range("E3").Value = 0.4832444488844488
$ws.Range("F3").Value = 0.6382978723404256
$ws.Range("G3").Value = 0.8956605459035781
$ws.Range("H3").Value = 0.5506599999999999
$ws.Range("I3").Value = 0.6521739130434783
$ws.Range("J3").Value = 0.7738100000000001
$ws.Range("K3").Value = 0.46872
$ws.Range("L3").Value = 0.625
$ws.Range("M3").Value = @'
[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0 0 1 0 1 1 0 1 1 0 1 1 1]
'@
$ws.Range("N3").Value = @'
[0 1 0 1 1 1 0 1 0 1 0 0 1 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 0 1 0 0 0 1 0 1]
'@
$ws.Range("O3").Value = 69
$ws.Rows.Item(3).AutoFit()

# Row 4
$ws.Range("A4").Value = @'
Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector', RandomUnderSampler(random_state=42)),
                ('model',
                 BaggingClassifier(estimator=DecisionTreeClassifier(class_weight='balanced',
                                                                    criterion='entropy',
                                                                    max_depth=2,
                                                                    min_samples_leaf=3,
                                                                    random_state=42),
                                   n_estimators=5, random_state=42))])
'@
$ws.Range("B4").Value = 0.6221428571428571
$ws.Range("C4").Value = @'
{'selector': RandomUnderSampler(random_state=42), 'scaler': MinMaxScaler(), 'model__n_estimators': 5, 'model__estimator__min_samples_split': 2, 'model__estimator__min_samples_leaf': 3, 'model__estimator__max_features': None, 'model__estimator__max_depth': 2, 'model__estimator__criterion': 'entropy', 'model__estimator__class_weight': 'balanced'}
'@
$ws.Range("D4").Value = 0.8608482332958938
$ws.Range("E4").Value = 0.5340914774114774
$ws.Range("F4").Value = 0.6274509803921569
$ws.Range("G4").Value = 0.8421631503872792
$ws.Range("H4").Value = 0.5173022222222222
$ws.Range("I4").Value = 0.64
$ws.Range("J4").Value = 0.8870842105263158
$ws.Range("K4").Value = 0.57896
$ws.Range("L4").Value = 0.6153846153846154
$ws.Range("M4").Value = @'
[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1 0 1 0 1 0 1 0 1 1 1 0 1]
'@
$ws.Range("N4").Value = @'
[1 1 1 1 0 1 0 1 1 1 1 0 0 0 1 0 1 0 1 1 1 1 1 0 1 1 1 1 1 1 1 0 1 0 1 0]
'@
$ws.Range("O4").Value = 23
$ws.Rows.Item(4).AutoFit()

# Row 5
$ws.Range("A5").Value = @'
Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector', RandomUnderSampler(random_state=42)),
                ('model',
                 BaggingClassifier(estimator=DecisionTreeClassifier(class_weight='balanced',
                                                                    max_depth=6,
                                                                    max_features='log2',
                                                                    min_samples_leaf=2,
                                                                    min_samples_split=4,
                                                                    random_state=42),
                                   n_estimators=5, random_state=42))])
'@
$ws.Range("B5").Value = 0.7016666666666667
$ws.Range("C5").Value = @'
{'selector': RandomUnderSampler(random_state=42), 'scaler': MinMaxScaler(), 'model__n_estimators': 5, 'model__estimator__min_samples_split': 4, 'model__estimator__min_samples_leaf': 2, 'model__estimator__max_features': 'log2', 'model__estimator__max_depth': 6, 'model__estimator__criterion': 'gini', 'model__estimator__class_weight': 'balanced'}
'@
$ws.Range("D5").Value = 0.8007655300297267
$ws.Range("E5").Value = 0.4523190342990343
$ws.Range("F5").Value = 0.5263157894736842
$ws.Range("G5").Value = 0.8844170814131247
$ws.Range("H5").Value = 0.5168366666666667
$ws.Range("I5").Value = 0.625
$ws.Range("J5").Value = 0.736
$ws.Range("K5").Value = 0.428
$ws.Range("L5").Value = 0.4545454545454545
$ws.Range("M5").Value = @'
[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0 1 1 1 1 1 0]
'@
$ws.Range("N5").Value = @'
[0 1 1 0 0 0 0 1 0 0 1 1 0 1 1 0 0 0 0 0 0 1 0 1 1 0 0 1 0 1 0 1 1 1 0 1]
'@
$ws.Range("O5").Value = 99
$ws.Rows.Item(5).AutoFit()

# Row 6
$ws.Range("A6").Value = @'
Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector', RandomUnderSampler(random_state=42)),
                ('model',
                 BaggingClassifier(estimator=DecisionTreeClassifier(max_depth=4,
                                                                    max_features='log2',
                                                                    random_state=42),
                                   n_estimators=5, random_state=42))])
'@
$ws.Range("B6").Value = 0.7566666666666666
$ws.Range("C6").Value = @'
{'selector': RandomUnderSampler(random_state=42), 'scaler': MinMaxScaler(), 'model__n_estimators': 5, 'model__estimator__min_samples_split': 2, 'model__estimator__min_samples_leaf': 1, 'model__estimator__max_features': 'log2', 'model__estimator__max_depth': 4, 'model__estimator__criterion': 'gini', 'model__estimator__class_weight': None}
'@
$ws.Range("D6").Value = 0.8106785276542967
$ws.Range("E6").Value = 0.5142481385281384
$ws.Range("F6").Value = 0.55
$ws.Range("G6").Value = 0.9054163495814522
$ws.Range("H6").Value = 0.5798438095238095
$ws.Range("I6").Value = 0.55
$ws.Range("J6").Value = 0.7370363636363636
$ws.Range("K6").Value = 0.49992
$ws.Range("L6").Value = 0.55
$ws.Range("M6").Value = @'
[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1 1 0 1 0 1 1 1 1 1 1 1 0]
'@
$ws.Range("N6").Value = @'
[1 1 1 0 1 1 0 1 0 1 0 0 1 1 1 0 0 1 0 1 0 1 1 1 0 0 1 0 1 0 1 1 0 1 0 0]
'@
$ws.Range("O6").Value = 89
$ws.Rows.Item(6).AutoFit()
